# Auto-generated script applying the Kujata_Profits market-data refresh
# (values only; no formulas are present in this workbook - every cell is a literal)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3537.8086
$ws.Range("I64").Value = 3581.2942
$ws.Range("J64").Value = 3513.1667
$ws.Range("K64").Value = 3581.2942
$ws.Range("L64").Value = 3513.1667
$ws.Range("M64").Value = -3333.2942
$ws.Range("N64").Value = -4009.1667
$ws.Range("H67").Value = 3537.8086
$ws.Range("I67").Value = 3581.2942
$ws.Range("J67").Value = 3513.1667
$ws.Range("K67").Value = 3581.2942
$ws.Range("L67").Value = 3513.1667
$ws.Range("M67").Value = -2723.2942
$ws.Range("N67").Value = -5229.1667
$ws.Range("H70").Value = 1755.5883
$ws.Range("I70").Value = 1767.3572
$ws.Range("J70").Value = 1700.6666
$ws.Range("K70").Value = 5302.071599999999
$ws.Range("L70").Value = 5101.9998
$ws.Range("M70").Value = -5032.071599999999
$ws.Range("N70").Value = -5641.9998
$ws.Range("H73").Value = 1755.5883
$ws.Range("I73").Value = 1767.3572
$ws.Range("J73").Value = 1700.6666
$ws.Range("K73").Value = 5302.071599999999
$ws.Range("L73").Value = 5101.9998
$ws.Range("M73").Value = -4366.071599999999
$ws.Range("N73").Value = -6973.9998
$ws.Range("H132").Value = 6540168.5
$ws.Range("I132").Value = 8774313
$ws.Range("J132").Value = 9592.691999999999
$ws.Range("K132").Value = 26322939
$ws.Range("L132").Value = 28778.076
$ws.Range("M132").Value = -26320409
$ws.Range("N132").Value = -33838.076
$ws.Range("H137").Value = 2899.2593
$ws.Range("I137").Value = 1886.6666
$ws.Range("J137").Value = 3405.5557
$ws.Range("K137").Value = 5659.9998
$ws.Range("L137").Value = 10216.6671
$ws.Range("M137").Value = -3109.9998
$ws.Range("N137").Value = -15316.6671
$ws.Range("H138").Value = 2042.45
$ws.Range("I138").Value = 805.1818
$ws.Range("J138").Value = 2195.3708
$ws.Range("K138").Value = 2415.5454
$ws.Range("L138").Value = 6586.1124
$ws.Range("M138").Value = 2724.4546
$ws.Range("N138").Value = -16866.1124

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3162.5098
$ws.Range("I32").Value = 3388.4187
$ws.Range("K32").Value = 3388.4187
$ws.Range("M32").Value = -3101.4187
$ws.Range("H63").Value = 55557596
$ws.Range("I63").Value = 2142.7856
$ws.Range("J63").Value = 250001680
$ws.Range("K63").Value = 2142.7856
$ws.Range("L63").Value = 250001680
$ws.Range("M63").Value = -1456.7856
$ws.Range("N63").Value = -250003052
$ws.Range("H66").Value = 55557596
$ws.Range("I66").Value = 2142.7856
$ws.Range("J66").Value = 250001680
$ws.Range("K66").Value = 10713.928
$ws.Range("L66").Value = 1250008400
$ws.Range("M66").Value = -7281.928
$ws.Range("N66").Value = -1250015264
$ws.Range("H88").Value = 2507.8572
$ws.Range("I88").Value = 1987.25
$ws.Range("J88").Value = 2716.1
$ws.Range("K88").Value = 1987.25
$ws.Range("L88").Value = 2716.1
$ws.Range("M88").Value = -1581.25
$ws.Range("N88").Value = -3528.1
$ws.Range("H91").Value = 2507.8572
$ws.Range("I91").Value = 1987.25
$ws.Range("J91").Value = 2716.1
$ws.Range("K91").Value = 1987.25
$ws.Range("L91").Value = 2716.1
$ws.Range("M91").Value = -583.25
$ws.Range("N91").Value = -5524.1
$ws.Range("H103").Value = 62833.168
$ws.Range("J103").Value = 62833.168
$ws.Range("L103").Value = 62833.168
$ws.Range("N103").Value = -65177.168
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = $null
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3638.8157
$ws.Range("I86").Value = 3862.276
$ws.Range("J86").Value = 2918.7778
$ws.Range("K86").Value = 3862.276
$ws.Range("L86").Value = 2918.7778
$ws.Range("M86").Value = -2739.276
$ws.Range("N86").Value = -5164.7778
$ws.Range("H89").Value = 3638.8157
$ws.Range("I89").Value = 3862.276
$ws.Range("J89").Value = 2918.7778
$ws.Range("K89").Value = 19311.38
$ws.Range("L89").Value = 14593.889
$ws.Range("M89").Value = -13695.38
$ws.Range("N89").Value = -25825.889
$ws.Range("H105").Value = 250003070
$ws.Range("I105").Value = 333336670
$ws.Range("K105").Value = 333336670
$ws.Range("M105").Value = -333334923

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1843.9
$ws.Range("I31").Value = 1218.9412
$ws.Range("J31").Value = 2165.8484
$ws.Range("K31").Value = 1218.9412
$ws.Range("L31").Value = 2165.8484
$ws.Range("M31").Value = -923.9412
$ws.Range("N31").Value = -2755.8484
$ws.Range("H34").Value = 1843.9
$ws.Range("I34").Value = 1218.9412
$ws.Range("J34").Value = 2165.8484
$ws.Range("K34").Value = 1218.9412
$ws.Range("L34").Value = 2165.8484
$ws.Range("M34").Value = -1016.9412
$ws.Range("N34").Value = -2569.8484
$ws.Range("H58").Value = 612.2857
$ws.Range("I58").Value = 547.7917
$ws.Range("J58").Value = 999.25
$ws.Range("K58").Value = 547.7917
$ws.Range("L58").Value = 999.25
$ws.Range("M58").Value = -344.7917
$ws.Range("N58").Value = -1405.25
$ws.Range("H132").Value = 6509.7393
$ws.Range("I132").Value = 8775.857
$ws.Range("J132").Value = 2984.6667
$ws.Range("K132").Value = 26327.571
$ws.Range("L132").Value = 8954.000100000001
$ws.Range("M132").Value = -23797.571
$ws.Range("N132").Value = -14014.0001
$ws.Range("H134").Value = 8773469
$ws.Range("I134").Value = 10418123
$ws.Range("K134").Value = 31254369
$ws.Range("M134").Value = -31251834
$ws.Range("H136").Value = 612.2857
$ws.Range("I136").Value = 547.7917
$ws.Range("J136").Value = 999.25
$ws.Range("K136").Value = 1643.3751
$ws.Range("L136").Value = 2997.75
$ws.Range("M136").Value = 906.6249
$ws.Range("N136").Value = -8097.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2854.17
$ws.Range("I68").Value = 971.2353000000001
$ws.Range("J68").Value = 3239.8313
$ws.Range("K68").Value = 2913.7059
$ws.Range("L68").Value = 9719.493899999999
$ws.Range("M68").Value = -2102.7059
$ws.Range("N68").Value = -11341.4939
$ws.Range("H71").Value = 2854.17
$ws.Range("I71").Value = 971.2353000000001
$ws.Range("J71").Value = 3239.8313
$ws.Range("K71").Value = 8741.117700000001
$ws.Range("L71").Value = 29158.4817
$ws.Range("M71").Value = -4685.117700000001
$ws.Range("N71").Value = -37270.4817

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 17979.666
$ws.Range("I136").Value = 20975.6
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 62926.8
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -60376.8
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = $null
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = $null
$ws.Range("H132").Value = 4350.5625
$ws.Range("I132").Value = 4144.88
$ws.Range("J132").Value = 5085.143
$ws.Range("K132").Value = 12434.64
$ws.Range("L132").Value = 15255.429
$ws.Range("M132").Value = -9904.639999999999
$ws.Range("N132").Value = -20315.429
